$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B values (parameter names) for the new rows 506-537
$bVals = @(
    "TIL_FIRSTDRUM_PSDT1",
    "TIL_FIRSTDRUM_PSDT2",
    "TIL_FIRSTDRUM_PSDT3",
    "TIL_FIRSTDRUM_PSDT4",
    "TIL_FIRSTDRUM_PSDT5",
    "TIL_FIRSTDRUM_PSDT6",
    "TIL_FIRSTDRUM_PSDT7",
    "TIL_FIRSTDRUM_PSDT8",
    "TIL_FIRSTDRUM_RESULT21",
    "TIL_FIRSTDRUM_RESULT22",
    "TIL_FIRSTDRUM_RESULT23",
    "TIL_FIRSTDRUM_RESULT24",
    "TIL_FIRSTDRUM_RESULT25",
    "TIL_FIRSTDRUM_RESULT26",
    "TIL_FIRSTDRUM_RESULT27",
    "TIL_FIRSTDRUM_RESULT28",
    "TIL_IPT_ERH_%1",
    "TIL_IPT_ERH_%2",
    "TIL_IPT_ERH_%3",
    "TIL_IPT_ERH_%4",
    "TIL_IPT_ERH_%5",
    "TIL_IPT_ERH_%6",
    "TIL_IPT_ERH_%7",
    "TIL_IPT_ERH_%8",
    "TIL_CF_SPRAY_QTY_LOT1",
    "TIL_CF_SPRAY_QTY_LOT2",
    "TIL_CF_SPRAY_QTY_LOT3",
    "TIL_CF_SPRAY_QTY_LOT4",
    "TIL_CF_SPRAY_QTY_LOT5",
    "TIL_CF_SPRAY_QTY_LOT6",
    "TIL_CF_SPRAY_QTY_LOT7",
    "TIL_CF_SPRAY_QTY_LOT8"
)

# Column E values (short labels) for the new rows 506-537
$eVals = @(
    "PSD1",
    "PSD2",
    "PSD3",
    "PSD4",
    "PSD5",
    "PSD6",
    "PSD7",
    "PSD8",
    "SV1",
    "SV2",
    "SV3",
    "SV4",
    "SV5",
    "SV6",
    "SV7",
    "SV8",
    "ERH1",
    "ERH2",
    "ERH3",
    "ERH4",
    "ERH5",
    "ERH6",
    "ERH7",
    "ERH8",
    "Totalspratamt1",
    "Totalspratamt2",
    "Totalspratamt3",
    "Totalspratamt4",
    "Totalspratamt5",
    "Totalspratamt6",
    "Totalspratamt7",
    "Totalspratamt8"
)

$startRow = 506
for ($i = 0; $i -lt 32; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "EMI_0501"
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
    $ws.Cells.Item($r, 3).Value = "AG-1749 CF GRANULES"
    $ws.Cells.Item($r, 4).Value = "CF"
    $ws.Cells.Item($r, 5).Value = $eVals[$i]
    $ws.Cells.Item($r, 6).Value = "tbd"
}

# Extend the hidden _FilterDatabase defined name to cover the newly added rows
$wb.Names("_xlnm._FilterDatabase").RefersTo = "=parameters!`$A`$1:`$F`$537"

# Leave the selection where the author last left it while adding these rows
[void]$ws.Range("E523").Select()
